# Generate Report for Handoff
# Replace the old file id "e005a30e-17c8-4e91-83d2-e3ed56a42def" with the new
# one "ebe840a9-4e14-4a32-93d4-fe2121ecce10" everywhere it is referenced, and
# refresh the associated handoff timestamps / xlf file names.

$wb = $excel.ActiveWorkbook

$oldId = "e005a30e-17c8-4e91-83d2-e3ed56a42def"
$newId = "ebe840a9-4e14-4a32-93d4-fe2121ecce10"

$oldZhXlf = "$oldId.41f0197e57ec089d05ca2e99709f34527ce5005f.zh-cn.xlf"
$newZhXlf = "$newId.a9ed4fc76c8f56ee302467c35ae58749087ab49f.zh-cn.xlf"

$oldDeXlf = "$oldId.41f0197e57ec089d05ca2e99709f34527ce5005f.de-de.xlf"
$newDeXlf = "$newId.a9ed4fc76c8f56ee302467c35ae58749087ab49f.de-de.xlf"

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item(1)

$wsOverview.Range("A2").Value = "$newId.md"
$wsOverview.Range("G2").Value = "2016-08-18 15:03:10"

# B2 holds a (external) hyperlink whose cached display text must be updated
# while keeping the same underlying target address.
$overviewLinkAddr = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d1afddaa9c107eec1528dd890c8d12976d84811d/e2e/$oldId.md"
$overviewNewDisplay = "e2e\$newId.md"
$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), $overviewLinkAddr, "", "", $overviewNewDisplay)

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item(2)

$wsZhCn.Range("G2").Value = $newZhXlf
$wsZhCn.Range("H2").Value = "2016-08-18 15:02:57"

$zhLinkAddr = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d1afddaa9c107eec1528dd890c8d12976d84811d/e2e/$oldId.md"
$zhNewDisplay = "$newId.md"
$wsZhCn.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), $zhLinkAddr, "", "", $zhNewDisplay)

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item(3)

$wsDeDe.Range("G2").Value = $newDeXlf
# H2 shares its value with Overview!G2 ("Latest HO Xliff Generate Date");
# both were regenerated together and now carry the later timestamp.
$wsDeDe.Range("H2").Value = "2016-08-18 15:03:10"

$deLinkAddr = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d1afddaa9c107eec1528dd890c8d12976d84811d/e2e/$oldId.md"
$deNewDisplay = "$newId.md"
$wsDeDe.Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), $deLinkAddr, "", "", $deNewDisplay)
